$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H), mirroring the formatting used by the
# other header cells (e.g. G1, which holds the "sum" header).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data values for the Save column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
